$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "official_cluster_311111111111111111111111111111111111111111"

$ws.Range("A2").Select()
